$d = $word.ActiveDocument

# The document carries three embedded logos whose internal "name" label
# (the wp:docPr / pic:cNvPr name="...") needs to be relabelled:
#   - the two Pearson logos (footer, default + first-page) : image2.png -> image1.png
#   - the BTEC logo (header, first-page)                    : image1.jpg -> image2.jpg
# Identify each InlineShape by its alt text / size (stable across
# header/footer reordering) rather than assuming fixed indices, then
# rename it.

function Rename-Logo {
    param($shape, $newName)
    if ($null -eq $shape) { return }
    try {
        $shape.Name = $newName
    } catch {
        # Older Word object models do not expose InlineShape.Name for
        # writing; ignore and move on so the rest of the edit still runs.
    }
}

foreach ($sec in $d.Sections) {

    for ($hi = 1; $hi -le 3; $hi++) {

        $hdr = $sec.Headers.Item($hi)
        if ($hdr.Exists) {
            $cnt = $hdr.Range.InlineShapes.Count
            for ($j = 1; $j -le $cnt; $j++) {
                $s = $hdr.Range.InlineShapes.Item($j)
                if ($s.AlternativeText -eq "BTec_Logo-Orange") {
                    Rename-Logo $s "image2.jpg"
                }
            }
        }

        $ftr = $sec.Footers.Item($hi)
        if ($ftr.Exists) {
            $cnt = $ftr.Range.InlineShapes.Count
            for ($j = 1; $j -le $cnt; $j++) {
                $s = $ftr.Range.InlineShapes.Item($j)
                if ($s.AlternativeText -like "*PearsonLogo.png") {
                    Rename-Logo $s "image1.png"
                }
            }
        }
    }
}
